$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.336.61'
$ws.Range('E2').Value = '  +1.39%  '
$ws.Range('D3').Value = '1.864.74'
$ws.Range('E3').Value = '  +1.33%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.64'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.44%  '
$ws.Range('E6').Value = '  +0.95%  '
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.77'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.75%  '
$ws.Range('E10').Value = '  +1.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0988'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.72%  '
$ws.Range('D12').Value = '2.134.21'
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.57'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.48%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.859.62'
$ws.Range('E14').Value = '  +1.03%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.682'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.39%  '
$ws.Range('E16').Value = '  +1.85%  '
$ws.Range('D17').Value = '35.340.29'
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.23'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E19').Value = '  +1.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '241.44'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.29'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.01%  '
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.24'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '169.56'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.19%  '
$ws.Range('E26').Value = '  +25.21%  '
$ws.Range('E27').Value = '  +5.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.77'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.81%  '
$ws.Range('E30').Value = '  +2.29%  '
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('E32').Value = '  +2.03%  '
$ws.Range('E33').Value = '  +27.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.05'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.66%  '
$ws.Range('E35').Value = '  +9.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.813'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +17.18%  '
$ws.Range('E37').Value = '  +6.89%  '
$ws.Range('E38').Value = '  +3.74%  '
$ws.Range('E39').Value = '  +4.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '90.85'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.46%  '
$ws.Range('D41').Value = '1.350.98'
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0605'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +15.56%  '
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '15.25'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.37'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.05'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +54.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.42'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.66%  '
$ws.Range('E47').Value = '  +6.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.73'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('D49').Value = '2.051.07'
$ws.Range('E49').Value = '  +1.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0687'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.36%  '
$ws.Range('E51').Value = '  -0.96%  '
